$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40..151 down to 41..152
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly data point
$ws.Cells.Item(40, 1).Value = 8
$ws.Cells.Item(40, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(40, 3).Value = "Coquimbo"
$ws.Cells.Item(40, 4).Value = 44742
$ws.Cells.Item(40, 5).Value = 4
$ws.Cells.Item(40, 6).Value = 100112044
$ws.Cells.Item(40, 7).Value = "Perejil"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 11).Value = 1300
$ws.Cells.Item(40, 12).Value = 1500
$ws.Cells.Item(40, 13).Value = 1400
$ws.Cells.Item(40, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(40, 16).Value = 933
$ws.Cells.Item(40, 17).Value = 1.5
$ws.Cells.Item(40, 18).Value = "Hortaliza"
